$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update odds values in row 4 (E9XIRWCr, Atl. Tucuman vs Central Cordoba) ---
$ws.Range("I4").Value = 4.5
$ws.Range("K4").Value = 2
$ws.Range("M4").Value = 1.08
$ws.Range("N4").Value = 8
$ws.Range("O4").Value = 1.44
$ws.Range("P4").Value = 2.63
$ws.Range("Q4").Value = 2.35
$ws.Range("R4").Value = 1.57
$ws.Range("S4").Value = 1.53
$ws.Range("T4").Value = 2.38
$ws.Range("U4").Value = 2.1
$ws.Range("V4").Value = 1.67
$ws.Range("Y4").Value = 9
$ws.Range("AA4").Value = 17
$ws.Range("AB4").Value = 34
$ws.Range("AC4").Value = 7.5
$ws.Range("AH4").Value = 9.5
$ws.Range("AJ4").Value = 15
$ws.Range("AO4").Value = 10
$ws.Range("AT4").Value = 2.38
$ws.Range("AV4").Value = 67
$ws.Range("BC4").Value = 401
$ws.Range("BD4").Value = 151

# --- Update odds values in row 5 (dn2dujsr, Ind. Rivadavia vs Dep. Riestra) ---
$ws.Range("G5").Value = 1.8
$ws.Range("I5").Value = 4.75
$ws.Range("N5").Value = 8
$ws.Range("O5").Value = 1.44
$ws.Range("P5").Value = 2.63
$ws.Range("AE5").Value = 21
$ws.Range("AF5").Value = 81
$ws.Range("AH5").Value = 10
$ws.Range("AP5").Value = 26
$ws.Range("AS5").Value = 251
$ws.Range("BD5").Value = 151

# --- Update odds values in row 7 (6k1eBeD4, The Strongest vs Oriente Petrolero) ---
$ws.Range("N7").Value = 34
$ws.Range("O7").Value = 1.07
$ws.Range("P7").Value = 9
$ws.Range("Q7").Value = 1.25
$ws.Range("R7").Value = 4
$ws.Range("U7").Value = 2
$ws.Range("V7").Value = 1.73
$ws.Range("W7").Value = 13
$ws.Range("X7").Value = 8
$ws.Range("Z7").Value = 7.5
$ws.Range("AC7").Value = 34
$ws.Range("AN7").Value = 3.5
$ws.Range("AQ7").Value = 8.5
$ws.Range("AR7").Value = 26
$ws.Range("AS7").Value = 81
$ws.Range("AU7").Value = 11
$ws.Range("AV7").Value = 51
$ws.Range("AX7").Value = 15
$ws.Range("AY7").Value = 51
$ws.Range("BB7").Value = 251
$ws.Range("BC7").Value = 351

# --- Remove row 19 (Rio Ave vs Vitoria Guimaraes) - all subsequent rows shift up ---
$ws.Rows.Item(19).Delete()
